# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The workbook records a single "last changed" date shared by every row;
# it advanced from 45175 (2023-09-06) to 45177 (2023-09-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 472 }

$range = $ws.Range("C2:C$lastRow")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    $current = $cell.Value2()
    if ($current -eq 45175) {
        $cell.Value2 = 45177
    }
}
